$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each data row in the table is followed by 3 blank rows; data lives in
# rows 1, 5, 9, 13, 17 (1-indexed), 5 columns each.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "54÷9=" },
    @{ Row = 1;  Col = 2; Text = "63÷2=" },
    @{ Row = 1;  Col = 3; Text = "90÷3=" },
    @{ Row = 1;  Col = 4; Text = "10÷7=" },
    @{ Row = 1;  Col = 5; Text = "11÷2=" },

    @{ Row = 5;  Col = 1; Text = "30÷9=" },
    @{ Row = 5;  Col = 2; Text = "47÷3=" },
    @{ Row = 5;  Col = 3; Text = "87÷5=" },
    @{ Row = 5;  Col = 4; Text = "29÷6=" },
    @{ Row = 5;  Col = 5; Text = "27÷7=" },

    @{ Row = 9;  Col = 1; Text = "63÷2=" },
    @{ Row = 9;  Col = 2; Text = "82÷3=" },
    @{ Row = 9;  Col = 3; Text = "33÷3=" },
    @{ Row = 9;  Col = 4; Text = "67÷8=" },
    @{ Row = 9;  Col = 5; Text = "42÷4=" },

    @{ Row = 13; Col = 1; Text = "82÷9=" },
    @{ Row = 13; Col = 2; Text = "84÷4=" },
    @{ Row = 13; Col = 3; Text = "36÷8=" },
    @{ Row = 13; Col = 4; Text = "57÷7=" },
    @{ Row = 13; Col = 5; Text = "73÷9=" },

    @{ Row = 17; Col = 1; Text = "80÷7=" },
    @{ Row = 17; Col = 2; Text = "17÷5=" },
    @{ Row = 17; Col = 3; Text = "80÷8=" },
    @{ Row = 17; Col = 4; Text = "70÷4=" },
    @{ Row = 17; Col = 5; Text = "94÷7=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
